$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value2 = "30.268.08"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value2 = "  -0.28%  "
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value2 = "1.928.75"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value2 = "  -0.71%  "
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value2 = "0.9988"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).Value2 = "  -0.15%  "
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value2 = "0.7460"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value2 = "  +2.77%  "
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value2 = "249.48"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value2 = "  -1.00%  "
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value2 = "0.9979"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value2 = "  -0.14%  "
$ws.Cells.Item(8,2).Value2 = "Cardano"
$ws.Cells.Item(8,3).Value2 = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value2 = "0.3233"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value2 = "  -3.60%  "
$ws.Cells.Item(9,2).Value2 = "Solana"
$ws.Cells.Item(9,3).Value2 = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value2 = "27.96"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value2 = "  -3.12%  "
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value2 = "0.07102"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value2 = "  -4.12%  "
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value2 = "0.7895"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value2 = "  -4.00%  "
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value2 = "0.08019"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value2 = "  -1.55%  "
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value2 = "1.928.95"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value2 = "  -0.68%  "
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value2 = "5.385"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value2 = "  -2.31%  "
$ws.Cells.Item(15,5).Value2 = "  -0.80%  "
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value2 = "14.63"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value2 = "  -1.92%  "
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value2 = "30.253.35"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value2 = "  -0.40%  "
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value2 = "254.83"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value2 = "  +0.19%  "
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value2 = "0.000008060"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value2 = "  -3.58%  "
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value2 = "5.749"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value2 = "  -2.22%  "
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value2 = "2.183.25"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value2 = "  -0.57%  "
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value2 = "0.9976"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value2 = "  -0.18%  "
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value2 = "1.001"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value2 = "  +0.10%  "
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value2 = "6.826"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value2 = "  -2.45%  "
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value2 = "9.576"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value2 = "  -3.70%  "
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value2 = "164.13"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value2 = "  +1.73%  "
$ws.Cells.Item(27,2).Value2 = "Stellar"
$ws.Cells.Item(27,3).Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value2 = "0.1350"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value2 = "  +2.24%  "
$ws.Cells.Item(28,2).Value2 = "EthereumClassic"
$ws.Cells.Item(28,3).Value2 = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value2 = "19.10"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value2 = "  -1.66%  "
$ws.Cells.Item(29,2).Value2 = "LidoDAOToken"
$ws.Cells.Item(29,3).Value2 = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value2 = "2.299"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value2 = "  -4.79%  "
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value2 = "1.356"
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value2 = "  +0.85%  "
$ws.Cells.Item(31,5).Value2 = "  -2.81%  "
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value2 = "4.439"
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value2 = "  -0.90%  "
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value2 = "4.150"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value2 = "  -2.81%  "
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value2 = "0.05115"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value2 = "  -3.94%  "
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value2 = "1.290"
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Value2 = "  -1.34%  "
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value2 = "0.7491"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value2 = "  -1.93%  "
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value2 = "2.764"
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value2 = "  +0.51%  "
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value2 = "0.01976"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value2 = "  -1.20%  "
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value2 = "2.798"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value2 = "  -1.72%  "
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value2 = "78.18"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value2 = "  -4.10%  "
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value2 = "6.401"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value2 = "  -3.36%  "
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value2 = "0.4515"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value2 = "  -1.48%  "
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value2 = "1.986"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value2 = "  -3.13%  "
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value2 = "0.8413"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value2 = "  -0.50%  "
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value2 = "0.9976"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value2 = "  -0.20%  "
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value2 = "101.58"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value2 = "  -1.52%  "
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value2 = "9.823"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value2 = "  -0.35%  "
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value2 = "7.527"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value2 = "  -0.11%  "
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value2 = "976.32"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value2 = "  +10.93%  "
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value2 = "36.79"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value2 = "  -1.02%  "
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value2 = "0.4203"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value2 = "  -0.68%  "
